# Updated cryptos list on Sat Feb 24 10:44:57 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$value)
    # Force the cell to keep a literal text value (some of these look like
    # plain numbers, e.g. "380.34", and would otherwise be auto-converted
    # to a numeric type by value-type inference). Restore the original
    # "Normal" cell style afterwards so no stray formatting is introduced.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "51.073.32"
$ws.Range("E2").Value = "  -0.02%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "2.959.42"
$ws.Range("E3").Value = "  +0.71%  "

# Row 4 - TetherUSD
Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  +0.08%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "380.34"
$ws.Range("E5").Value = "  +1.11%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "102.24"
$ws.Range("E6").Value = "  -0.15%  "

# Row 7 - XRP
Set-TextValue $ws.Range("D7") "0.545"
$ws.Range("E7").Value = "  +1.79%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.01%  "

# Row 9 - Cardano
Set-TextValue $ws.Range("D9") "0.588"
$ws.Range("E9").Value = "  +0.88%  "

# Row 10 - Avalanche
Set-TextValue $ws.Range("D10") "36.48"
$ws.Range("E10").Value = "  -0.49%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  -0.54%  "

# Row 12 - Dogecoin
Set-TextValue $ws.Range("D12") "0.0851"
$ws.Range("E12").Value = "  +1.91%  "

# Row 13 - now WrappedliquidstakedEther2.0 (was Chainlink)
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue $ws.Range("D13") "3.427.00"
$ws.Range("E13").Value = "  +0.76%  "

# Row 14 - now Chainlink (was WrappedliquidstakedEther2.0)
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Range("D14") "18.41"
$ws.Range("E14").Value = "  +2.59%  "

# Row 15 - Uniswap
Set-TextValue $ws.Range("D15") "12.38"
$ws.Range("E15").Value = "  +73.51%  "

# Row 16 - Polkadot
Set-TextValue $ws.Range("D16") "7.74"
$ws.Range("E16").Value = "  +5.44%  "

# Row 17 - WrappedEther
Set-TextValue $ws.Range("D17") "2.961.68"
$ws.Range("E17").Value = "  +0.76%  "

# Row 18 - Polygon
$ws.Range("E18").Value = "  +3.82%  "

# Row 19 - WrappedBTC
Set-TextValue $ws.Range("D19") "51.128.00"
$ws.Range("E19").Value = "  +0.18%  "

# Row 20 - ImmutableX
$ws.Range("E20").Value = "  -2.06%  "

# Row 21 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D21") "12.37"
$ws.Range("E21").Value = "  -1.26%  "

# Row 22 - ShibaInu
$ws.Range("E22").Value = "  +0.78%  "

# Row 23 - PancakeSwap
Set-TextValue $ws.Range("D23") "3.34"
$ws.Range("E23").Value = "  +16.69%  "

# Row 24
Set-TextValue $ws.Range("D24") "269.76"
$ws.Range("E24").Value = "  +2.62%  "

# Row 25
Set-TextValue $ws.Range("D25") "69.72"
$ws.Range("E25").Value = "  +2.25%  "

# Row 26
Set-TextValue $ws.Range("D26") "7.95"
$ws.Range("E26").Value = "  -2.34%  "

# Row 27
$ws.Range("E27").Value = "  +0.06%  "

# Row 28 - Kaspa
$ws.Range("E28").Value = "  -0.12%  "

# Row 29 - EthereumClassic
Set-TextValue $ws.Range("D29") "25.85"
$ws.Range("E29").Value = "  +0.75%  "

# Row 30 - RenderToken
$ws.Range("E30").Value = "  -9.91%  "

# Row 31 - Hedera
$ws.Range("E31").Value = "  -3.56%  "

# Row 32 - Cosmos
Set-TextValue $ws.Range("D32") "10.54"
$ws.Range("E32").Value = "  +7.03%  "

# Row 33 - OKB
Set-TextValue $ws.Range("D33") "51.16"
$ws.Range("E33").Value = "  +0.96%  "

# Row 34 - InjectiveProtocol
Set-TextValue $ws.Range("D34") "34.22"
$ws.Range("E34").Value = "  +0.43%  "

# Row 35 - Toncoin
$ws.Range("E35").Value = "  +2.07%  "

# Row 36 - VeChain
$ws.Range("E36").Value = "  -3.90%  "

# Row 37 - FirstDigitalUSD
$ws.Range("E37").Value = "  +0.01%  "

# Row 38 - LidoDAOToken
Set-TextValue $ws.Range("D38") "3.26"
$ws.Range("E38").Value = "  +9.50%  "

# Row 39 - Stellar
$ws.Range("E39").Value = "  +2.29%  "

# Row 40 - Celestia
Set-TextValue $ws.Range("D40") "16.71"
$ws.Range("E40").Value = "  +1.48%  "

# Row 41 - ARBITRUM
$ws.Range("E41").Value = "  +2.91%  "

# Row 42 - Stacks
Set-TextValue $ws.Range("D42") "2.50"
$ws.Range("E42").Value = "  -2.58%  "

# Row 43 - Monero
Set-TextValue $ws.Range("D43") "123.92"
$ws.Range("E43").Value = "  +1.90%  "

# Row 44 - EnergySwap
Set-TextValue $ws.Range("D44") "21.78"
$ws.Range("E44").Value = "  +3.20%  "

# Row 45 - NEARProtocol
Set-TextValue $ws.Range("D45") "3.56"
$ws.Range("E45").Value = "  +10.54%  "

# Row 46 - Maker
Set-TextValue $ws.Range("D46") "2.090.33"
$ws.Range("E46").Value = "  +4.38%  "

# Row 47 - WEMIXToken
Set-TextValue $ws.Range("D47") "2.04"
$ws.Range("E47").Value = "  -0.57%  "

# Row 48 - ApeXProtocol
$ws.Range("E48").Value = "  -0.53%  "

# Row 49 - TheGraph
Set-TextValue $ws.Range("D49") "0.262"
$ws.Range("E49").Value = "  -3.23%  "

# Row 50 - BEAM
Set-TextValue $ws.Range("D50") "0.0321"
$ws.Range("E50").Value = "  -7.07%  "

# Row 51 - THORChain
Set-TextValue $ws.Range("D51") "5.36"
$ws.Range("E51").Value = "  +6.68%  "
